$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fax_Geschäftlich" header lived in column N (14). Remove that whole
# column so everything to its right shifts left: E-Mail_Geschäftlich
# (was O) -> N, Internetadresse (was P) -> O, Bemerkung (was Q) -> P, and
# the now-unused trailing column Q disappears.
$ws.Columns.Item(14).Delete()

# After the delete, Excel no longer flags column K (11) as "best fit" -
# re-apply its width explicitly (without AutoFit) so the stored width
# stays 10 characters but the bestFit flag is dropped, matching the
# post-edit workbook.
$ws.Columns.Item(11).ColumnWidth = 9.14

# Re-apply the AutoFilter over the new, narrower header range (A2:P2
# instead of A2:Q2).
$ws.AutoFilterMode = $false
$ws.Range("A2:P2").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the autofilter range.
$wb.Names.Item("Tabelle1!_FilterDatabase").RefersTo = "=Tabelle1!`$A`$2:`$P`$2"

# Leave the selection where the editor's session ended.
$ws.Range("L5").Select()
